$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "67.370.69"
$ws.Range("E2").Value = "  -2.22%  "
Set-TextValue "D3" "2.644.04"
$ws.Range("E3").Value = "  -3.27%  "
Set-TextValue "D5" "598.68"
$ws.Range("E5").Value = "  -0.92%  "
Set-TextValue "D6" "166.33"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E8").Value = "  -0.78%  "
Set-TextValue "D9" "2.643.68"
$ws.Range("E9").Value = "  -3.21%  "
Set-TextValue "D10" "0.145"
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("E13").Value = "  -2.14%  "
Set-TextValue "D14" "28.02"
$ws.Range("E14").Value = "  -2.28%  "
Set-TextValue "D15" "3.124.27"
$ws.Range("E16").Value = "  -3.16%  "
Set-TextValue "D17" "67.347.82"
$ws.Range("E17").Value = "  -2.00%  "
Set-TextValue "D18" "2.631.93"
$ws.Range("E18").Value = "  -3.88%  "
Set-TextValue "D19" "11.91"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("E20").Value = "  +2.53%  "
Set-TextValue "D21" "364.29"
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("E23").Value = "  -3.41%  "
Set-TextValue "D24" "11.14"
$ws.Range("E24").Value = "  +11.61%  "
Set-TextValue "D25" "2.02"
$ws.Range("E25").Value = "  -6.01%  "
$ws.Range("E26").Value = "  -0.05%  "
Set-TextValue "D27" "71.02"
$ws.Range("E27").Value = "  -3.66%  "
Set-TextValue "D28" "2.776.75"
$ws.Range("E28").Value = "  -3.47%  "
$ws.Range("E29").Value = "  -3.83%  "
$ws.Range("E30").Value = "  +0.46%  "
Set-TextValue "D31" "555.24"
$ws.Range("E31").Value = "  -6.19%  "
Set-TextValue "D32" "8.05"
$ws.Range("E32").Value = "  -2.91%  "
$ws.Range("E33").Value = "  -4.06%  "
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -5.20%  "
Set-TextValue "D38" "157.84"
$ws.Range("E38").Value = "  -2.05%  "
Set-TextValue "D39" "19.43"
$ws.Range("E39").Value = "  -2.49%  "
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("E41").Value = "  -4.17%  "
Set-TextValue "D42" "1.83"
$ws.Range("E42").Value = "  -5.17%  "
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("E45").Value = "  +0.07%  "
Set-TextValue "D46" "40.16"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("E47").Value = "  -3.45%  "
Set-TextValue "D48" "0.597"
$ws.Range("E48").Value = "  -1.55%  "
Set-TextValue "D49" "154.21"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("E50").Value = "  -2.24%  "
Set-TextValue "D51" "1.74"
$ws.Range("E51").Value = "  -3.39%  "
